# Swap the deck's theme palette from "Integral" to the stock "Office Theme"
# palette (the author's edit rewrote ppt/theme/theme1.xml <-> theme2.xml so
# the slide-master theme became "Office Theme" and the notes-master theme
# became "Integral"). The PowerPoint object model only exposes the single
# active theme color scheme used by the slide master, so we drive that via
# ThemeColorScheme, writing all twelve theme colors to the "Office" values.

function RGBVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Order matches the OOXML <a:clrScheme> child order / ThemeColorScheme index:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $cs.Item($i).RGB = RGBVal($officeColors[$i - 1])
}
